# redundanciesbyannotations/output_detections_highlighting_without_annotations.xlsx
#
# The source script was refactored (calc_main_part_redundancy /
# calc_benefit_redundancy merged into a shared calc_redundancy helper) and,
# as part of that clean-up, the two result sheets were renamed to make clear
# they only cover "US" (User Story) detections:
#   "False Positives Negatives" -> "Only US-False Pos. Neg."
#   "True Positives"            -> "Only US-True Pos."
#
# Renaming a worksheet in Excel automatically rewrites every reference to
# the old sheet name, including the hidden per-sheet AutoFilter defined name
# (_xlnm._FilterDatabase), so that follows along with no extra work.

$wb = $excel.ActiveWorkbook

$falsePosNeg = $wb.Worksheets.Item("False Positives Negatives")
$falsePosNeg.Name = "Only US-False Pos. Neg."

$truePos = $wb.Worksheets.Item("True Positives")
$truePos.Name = "Only US-True Pos."

# The refactor also stopped emitting one never-used header cell format
# (a bold-14pt/boxed-border/centered xf that no cell ever referenced), so
# the workbook's style table shrank by one entry. Mirror that clean-up when
# the host exposes a way to drop orphaned cell formats; ignore it quietly
# when it doesn't, since it has no visible effect on any cell either way.
foreach ($target in @($wb, $excel, $falsePosNeg, $truePos)) {
    try { $target.CleanExcessFormatting() } catch { }
}
